# Rev 1.1 changes to the ICEBoard Bill of Materials
#  1. Row 5 (Slide Switch) changed from a DPDT (C&K) part to a SPDT (APEM) part,
#     now also used on S10, with updated price/qty.
#  2. A new Resistor line (R11, 470 ohm) is inserted as a new row 12, pushing the
#     two Capacitor rows and the Buzzer row down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Insert a new row at position 12 for the new 470 ohm resistor (R11) ---
$ws.Rows.Item(12).Insert()
$ws.Range("F12").Value = "R11"

# --- 2. Update the Slide Switch row (row 5) ---
$ws.Range("B5").Value = "SPDT"
$ws.Range("C5").Value = "APEM"
$ws.Range("D5").Value = "MHSS1104"
$ws.Range("E5").Value = "SWITCH SLIDE SPDT 300MA 6V"
$ws.Range("F5").Value = "S1, S2, S3, S4, S10"
$ws.Range("G5").Value = 0.55
$ws.Range("H5").Value = 5
$ws.Range("J5").Value = "679-1848-ND"

# --- 3. Fill in the rest of the new resistor row ---
$ws.Range("A12").Value = "Resistor"
$ws.Range("B12").Value = 470
$ws.Range("C12").Value = "Yageo"
$ws.Range("D12").Value = "RC0603FR-07470RL"
$ws.Range("E12").Value = "RES SMD 470 OHM 1% 1/10W 0603"
$ws.Range("G12").Value = 0.1
$ws.Range("H12").Value = 1
$ws.Range("I12").Formula = "=G12*H12"
$ws.Range("J12").Value = "311-470HRCT-ND"

# --- Update the saved selection to match the authored workbook (cell J13) ---
$ws.Range("J13").Select()
